$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for B1 (Active cases vs Activecases)
$ws.Range("B1").Value = "Active cases"

# New cluster names and active-case counts (rows 2-53), replacing the old 54-row table
$names = @(
    'Al Iman College Melton South',
    'Al Siraat College Epping',
    'Ashwood High School Ashwood',
    'Bairnsdale West Primary School',
    'Baptcare Westhaven community',
    'Belvedere Aged Care Noble Park CLOSED',
    'Bethany Catholic Primary School Werribee',
    'Black Rock Primary School Black Rock',
    'Blue Cross The Gables Camberwell',
    'BlueCross Elly Kay Mordialloc',
    'Brookside P-9 College Caroline Springs',
    'Bupa Aged Care Eastwood',
    'Camp Coolamatong Farm Camp Banksia Peninsula',
    'Christ the Priest Primary School Caroline Springs',
    'Covenant College Bell Post Hill',
    'Creekside K-9 College Caroline Springs',
    'Dandenong North Primary School Dandenong',
    'Epping Views Primary School Camp Cape Schanck',
    'Epping Views Primary School Epping',
    'Fitzroy Primary School Fitzroy',
    'Fronditha Thalpori St Albans Aged Care',
    'Hamilton Country Music Festival Hamilton Golf Club Hamilton',
    'Hodges Real Estate Brighton',
    'Holy Rosary Primary School White Hills',
    'Hope Aged Care Sunshine West',
    'Islamic College of Melbourne Tarneit Oct Nov',
    'Kerala Manor Aged Care Diamond Creek',
    'Little Munchkins Childcare Centre Hillside Exposure Site',
    'Lucknow Primary School Bairnsdale',
    'Nar Nar Goon Primary School Nar Nar Goon',
    'Nazareth Catholic Primary School Grovedale',
    'Nhill College Nhill',
    'Oakleigh Grammar Melbourne Private School Oakleigh',
    'Our Lady of the Southern Cross Primary School Manor Lakes',
    'Our Lady''s Catholic Primary School Craigieburn',
    'Parkdale Primary School Parkdale',
    'Parktone Primary School Parkdale',
    'Rochester and Elmore District Health Service Yalunkan Aged Care Hostel Rochester',
    'Social Gathering 20 November Sunbury',
    'Springside Primary School Caroline Springs Nov',
    'St Anne''s Catholic Primary School Sunbury',
    'St Brendan''s Primary School Somerville',
    'St Brendans Primary School Shepparton',
    'St Josephs Catholic Primary School Warragul',
    'St Kevin''s College Toorak Glendalough Campus Junior School',
    'St Mary of the Cross MacKillop Primary School Epping',
    'St Mary''s School Mooroopna',
    'Stawell Primary School',
    'Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North',
    'Village Glen Aged Care Residences Mornington',
    'Wagstaff Meat Processing Plant Cranbourne East',
    'Willmott Park Primary School Craigieburn'
)

$values = @(
    10,
    31,
    22,
    12,
    21,
    23,
    17,
    17,
    31,
    14,
    22,
    15,
    14,
    27,
    12,
    10,
    17,
    15,
    17,
    19,
    23,
    20,
    12,
    47,
    18,
    15,
    10,
    11,
    19,
    18,
    30,
    41,
    21,
    44,
    16,
    12,
    27,
    10,
    14,
    21,
    12,
    15,
    17,
    14,
    16,
    10,
    14,
    22,
    14,
    10,
    37,
    10
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove the two now-unused trailing rows (table shrank from 55 to 53 rows)
$ws.Rows.Item(54).Delete()
$ws.Rows.Item(54).Delete()
